# Generate Report for Handback
#
# The localization status report is updated to reflect that the files have
# been handed back (rather than still being "Ready for handoff"), and the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns are populated for each row on the per-language sheets.

$wb = $excel.ActiveWorkbook

function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime,
        [string]$Row2SourceUrl,
        [string]$Row2TargetUrl,
        [string]$Row3SourceUrl,
        [string]$Row3TargetUrl
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("C2").Value2 = "Handed back: in sync with en-US"
    $ws.Range("C3").Value2 = "Handed back: in sync with en-US"

    $row2Source = $ws.Range("A2").Value2
    $row2Target = $ws.Range("D2").Value2
    $row3Source = $ws.Range("A3").Value2
    $row3Target = $ws.Range("D3").Value2

    # Latest Target File (F) / Latest Handback File (G) now mirror the
    # source (.md) / target (.xlf) files that were handed off, since the
    # handback is in sync with them.
    $ws.Range("F2").Value2 = $row2Source
    $ws.Range("G2").Value2 = $row2Target
    $ws.Range("F3").Value2 = $row3Source
    $ws.Range("G3").Value2 = $row3Target

    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("G2").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
    $ws.Range("G3").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("F2"), $Row2SourceUrl, [Type]::Missing, [Type]::Missing, $row2Source) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $Row2TargetUrl, [Type]::Missing, [Type]::Missing, $row2Target) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $Row3SourceUrl, [Type]::Missing, [Type]::Missing, $row3Source) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $Row3TargetUrl, [Type]::Missing, [Type]::Missing, $row3Target) | Out-Null

    # Latest Handback DateTime (H) now records when the handback completed.
    $ws.Range("H2").Value2 = $HandbackDateTime
    $ws.Range("H3").Value2 = $HandbackDateTime
}

Update-LanguageSheet "zh-cn" "2016-03-13 00:41:58" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d8d3359bb86ea84e5d9e6cd1ec8649ac59834430/e2e/54b515d0-d368-4f1f-99f3-4f8cf3c71886.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd30e6b7d51e0c284e2299dd023529d40b29c2e2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/54b515d0-d368-4f1f-99f3-4f8cf3c71886.148dd6638baacae9ea6e78a30395944c9793c711.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d8d3359bb86ea84e5d9e6cd1ec8649ac59834430/e2e/b47a0c2b-7634-43d8-91d5-7a00f8ae86fc.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd30e6b7d51e0c284e2299dd023529d40b29c2e2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b47a0c2b-7634-43d8-91d5-7a00f8ae86fc.a9f1b68ee116756696689e3535fe3707bd8da256.zh-cn.xlf"

Update-LanguageSheet "de-de" "2016-03-13 00:42:04" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d8d3359bb86ea84e5d9e6cd1ec8649ac59834430/e2e/54b515d0-d368-4f1f-99f3-4f8cf3c71886.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/98b35a91bddcba8e29929c9230e19467379b733e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/54b515d0-d368-4f1f-99f3-4f8cf3c71886.148dd6638baacae9ea6e78a30395944c9793c711.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d8d3359bb86ea84e5d9e6cd1ec8649ac59834430/e2e/b47a0c2b-7634-43d8-91d5-7a00f8ae86fc.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/98b35a91bddcba8e29929c9230e19467379b733e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b47a0c2b-7634-43d8-91d5-7a00f8ae86fc.a9f1b68ee116756696689e3535fe3707bd8da256.de-de.xlf"
